# Append: 2025-09-24 01:43 JST
# Update the "取得日時" (acquired datetime) column (A) for all existing
# data rows on the active sheet ("ランサーズ") from the old timestamp
# to the new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-24 01:43:34"

# Data rows are 2 through 15 (row 1 is the header "取得日時" etc.)
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
